# Add "Admin - Manage Users - Use Case - 01" section.
#
# Structurally this inserts 4 new rows above the old row 6 ("2.2 System
# Requirements" row), which pushes that row (and everything below it) down
# by 4. The two new content rows become:
#   row 6: A6 = checkmark, D6 = "  -Overview"
#   row 7: D7 = "  -Admin", E7 = "Manage Users"
# rows 8-9 stay blank (just the checkmark-column styling carried down by
# the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 6; Excel copies row 5's formatting down
# into all 4 new rows (for every column that had a styled cell in row 5).
$ws.Rows("6:9").Insert()

# Row 5 had styled cells in columns A, C, D, E. The new section only wants
# A (checkmark column) and D/E (text), so drop the stray C cells (all 4
# rows) and the E cells that aren't actually part of the new content
# (rows 6, 8, 9 -- row 7 keeps its E cell).
$ws.Range("C6:C9").Clear()
$ws.Range("E6").Clear()
$ws.Range("E8").Clear()
$ws.Range("E9").Clear()
$ws.Range("D8").Clear()
$ws.Range("D9").Clear()

# Fill in the new content.
$ws.Range("A6").Value = "√"
$ws.Range("D6").Value = "  -Overview"
$ws.Range("D7").Value = "  -Admin"
$ws.Range("E7").Value = "Manage Users"

# Match the saved selection state from the target workbook.
$ws.Range("I6").Select()
